$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate header row to English
$ws.Range("A1").Value = "Age group"
$ws.Range("B1").Value = "Year-Month"
$ws.Range("C1").Value = "Deaths"

# Remove the now-unneeded "fecha" date column (column D), including its
# per-cell date formatting / number format
$ws.Columns.Item(4).Delete()
